$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# Reformat the multi-line remark text in column G ("備註說明") for the
# BatchType / Status / ErrCode rows: drop the dash-before-colon style and
# use a plain colon separator instead.
$ws.Range("G13").Value = "D:日批" + [char]10 + "M:月批"
$ws.Range("G14").Value = "S:成功" + [char]10 + "F:失敗"
$ws.Range("G15").Value = "五碼:一般交易錯誤" + [char]10 + "DB000:DB異常" + [char]10 + "LG000:邏輯錯誤"

# Leave the selection where the edit session ended up.
$ws.Range("G18").Select() | Out-Null
